$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated voltage magnitude (vm_pu) results for case with 380 kV
$data = @{
  2 = @{ "B"="1.02"; "C"="1.039015814634125"; "D"="1.050676651512808"; "E"="1.047706388479669"; "F"="1.058835101642867"; "I"="1.041474213513059"; "J"="1.04411011898795"; "K"="1.053430266320456"; "L"="1.050468276395659"; "M"="1.06156625132039"; "N"="1.045592875118017" }
  3 = @{ "B"="1.02"; "C"="1.03987993895287"; "D"="1.05127458493067"; "E"="1.048453102368312"; "F"="1.05961274713781"; "I"="1.041646872184187"; "J"="1.044619702400602"; "K"="1.053841246732026"; "L"="1.051027057440351"; "M"="1.062158100284983"; "N"="1.046103182197565" }
  4 = @{ "B"="1.02"; "C"="1.040439834849875"; "D"="1.051662024963069"; "E"="1.048937303034602"; "F"="1.060116980632773"; "I"="1.041757730985709"; "J"="1.044949526346278"; "K"="1.05410701673769"; "L"="1.051388984944139"; "M"="1.062541447857345"; "N"="1.046433474531076" }
  5 = @{ "B"="1.02"; "C"="1.040675392926385"; "D"="1.051825031442816"; "E"="1.049141104935585"; "F"="1.06032920849636"; "I"="1.041804128711337"; "J"="1.045088204424576"; "K"="1.054218706341968"; "L"="1.05154122368849"; "M"="1.062702696906397"; "N"="1.046572349548151" }
  6 = @{ "B"="1.02"; "C"="1.040714954545236"; "D"="1.051852408321907"; "E"="1.049175338421403"; "F"="1.060364856978024"; "I"="1.041811906915311"; "J"="1.045111490224193"; "K"="1.054237457126284"; "L"="1.051566790151151"; "M"="1.062729776544948"; "N"="1.046595668416275" }
  7 = @{ "B"="1.02"; "C"="1.040442981691666"; "D"="1.051664202568216"; "E"="1.048940025289431"; "F"="1.060119815459565"; "I"="1.041758351769908"; "J"="1.044951379292245"; "K"="1.05410850929896"; "L"="1.051391018835825"; "M"="1.062543602124531"; "N"="1.046435330108438" }
  8 = @{ "B"="1.02"; "C"="1.039307693283767"; "D"="1.050878613553734"; "E"="1.047958530183186"; "F"="1.059097693007857"; "I"="1.041532742547837"; "J"="1.044282315787802"; "K"="1.053569191753245"; "L"="1.050657043885648"; "M"="1.061766189417765"; "N"="1.045765316457077" }
  9 = @{ "B"="1.02"; "C"="1.037312984092351"; "D"="1.049498503042486"; "E"="1.04623695837556"; "F"="1.057304675225325"; "I"="1.041128613965311"; "J"="1.043104083169851"; "K"="1.052617666721462"; "L"="1.049366499560059"; "M"="1.060399285550782"; "N"="1.044585410613663" }
  10 = @{ "B"="1.02"; "C"="1.035987184431184"; "D"="1.048581365186581"; "E"="1.045094700457768"; "F"="1.056114892918572"; "I"="1.040854817290938"; "J"="1.04231917245213"; "K"="1.051982606176102"; "L"="1.048508116370531"; "M"="1.059490129921601"; "N"="1.043799385232692" }
  11 = @{ "B"="1.02"; "C"="1.035414069319703"; "D"="1.048184954770764"; "E"="1.044601408293012"; "F"="1.055601049301829"; "I"="1.04073523127238"; "J"="1.041979450994857"; "K"="1.051707466759326"; "L"="1.048136915171094"; "M"="1.059096977732995"; "N"="1.043459181331983" }
  12 = @{ "B"="1.02"; "C"="1.035201335458307"; "D"="1.048037819871397"; "E"="1.044418376833388"; "F"="1.05541038832099"; "I"="1.040690657653626"; "J"="1.041853287195363"; "K"="1.051605245991994"; "L"="1.047999108799537"; "M"="1.058951023053354"; "N"="1.043332838365423" }
  13 = @{ "B"="1.02"; "C"="1.035246960932897"; "D"="1.048069375816667"; "E"="1.044457628656749"; "F"="1.055451276509812"; "I"="1.040700225805259"; "J"="1.041880348656826"; "K"="1.051627173654764"; "L"="1.04802866534842"; "M"="1.058982327200661"; "N"="1.043359938257265" }
  14 = @{ "B"="1.02"; "C"="1.035396481672049"; "D"="1.048172790307385"; "E"="1.044586274780323"; "F"="1.055585285043879"; "I"="1.040731549941392"; "J"="1.041969021761077"; "K"="1.051699017594148"; "L"="1.048125522529389"; "M"="1.059084911441086"; "N"="1.043448737287494" }
  15 = @{ "B"="1.02"; "C"="1.035488625810642"; "D"="1.048236521967069"; "E"="1.044665564358821"; "F"="1.055667879142149"; "I"="1.04075082938485"; "J"="1.042023659390779"; "K"="1.051743280177825"; "L"="1.048185209323416"; "M"="1.059148127568919"; "N"="1.043503452508898" }
  16 = @{ "B"="1.02"; "C"="1.036025240663078"; "D"="1.04860768889985"; "E"="1.045127466464503"; "F"="1.056149023438233"; "I"="1.040862732177774"; "J"="1.042341721921618"; "K"="1.052000863147835"; "L"="1.048532762105813"; "M"="1.059516233200081"; "N"="1.043821966725013" }
  17 = @{ "B"="1.02"; "C"="1.036362104471371"; "D"="1.048840705276168"; "E"="1.045417558365396"; "F"="1.056451192776492"; "I"="1.040932650535282"; "J"="1.04254127525931"; "K"="1.052162397896107"; "L"="1.048750903532632"; "M"="1.059747276070455"; "N"="1.044021803451334" }
  18 = @{ "B"="1.02"; "C"="1.036558684418867"; "D"="1.048976688770232"; "E"="1.045586890623412"; "F"="1.056627572186372"; "I"="1.040973333276432"; "J"="1.042657685734651"; "K"="1.052256603428209"; "L"="1.048878188371345"; "M"="1.059882089352818"; "N"="1.044138379242903" }
  19 = @{ "B"="1.02"; "C"="1.036625728823177"; "D"="1.049023067268821"; "E"="1.045644649955126"; "F"="1.056687734848731"; "I"="1.040987188148337"; "J"="1.042697381098692"; "K"="1.05228872250114"; "L"="1.048921597076251"; "M"="1.059928065608005"; "N"="1.044178130978914" }
  20 = @{ "B"="1.02"; "C"="1.036325952542797"; "D"="1.048815697674757"; "E"="1.045386421145104"; "F"="1.056418759492746"; "I"="1.040925159238892"; "J"="1.042519863576743"; "K"="1.052145068282987"; "L"="1.048727494187438"; "M"="1.059722482184967"; "N"="1.044000361361722" }
  21 = @{ "B"="1.02"; "C"="1.035352447480718"; "D"="1.04814233427677"; "E"="1.044548386188467"; "F"="1.055545817238491"; "I"="1.040722330010557"; "J"="1.041942909072215"; "K"="1.051677861933475"; "L"="1.048096998441048"; "M"="1.059054700729597"; "N"="1.043422587515619" }
  22 = @{ "B"="1.02"; "C"="1.034741215261692"; "D"="1.047719598941383"; "E"="1.044022633902351"; "F"="1.054998141744917"; "I"="1.040593912322523"; "J"="1.041580293673363"; "K"="1.051383985685175"; "L"="1.047701011257956"; "M"="1.058635301104669"; "N"="1.043059457161309" }
  23 = @{ "B"="1.02"; "C"="1.035065159995107"; "D"="1.047943638101192"; "E"="1.044301235144904"; "F"="1.05528836252929"; "I"="1.040662073169441"; "J"="1.041772509356325"; "K"="1.051539786432704"; "L"="1.047910890244518"; "M"="1.058857588510384"; "N"="1.043251945812589" }
  24 = @{ "B"="1.02"; "C"="1.036342287756059"; "D"="1.048826997320507"; "E"="1.045400490323873"; "F"="1.056433414298831"; "I"="1.040928544540661"; "J"="1.042529538550145"; "K"="1.052152898835598"; "L"="1.048738071718986"; "M"="1.059733685320662"; "N"="1.044010050074696" }
  25 = @{ "B"="1.02"; "C"="1.03782796521437"; "D"="1.049854785531768"; "E"="1.046681072445271"; "F"="1.057767241723818"; "I"="1.041233865698154"; "J"="1.043408588457372"; "K"="1.052863789978759"; "L"="1.04969979372343"; "M"="1.060752297659661"; "N"="1.04489034833362" }
}

foreach ($rowKey in $data.Keys) {
  $rowNum = [int]$rowKey
  $rowData = $data[$rowKey]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$rowNum").Value = [double]$rowData[$col]
  }
}
